$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the four "set length" columns (B:E) on row 1 (submax reps header values)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 data (lichtwark passive values removed/replaced)
$ws.Range("B2").Value = 20.995868953125068
$ws.Range("C2").Value = 21.041776717499999
$ws.Range("D2").Value = 28.206329338124988
$ws.Range("E2").Value = 29.825140117499984

# Row 3 data
$ws.Range("B3").Value = 18.391727160000016
$ws.Range("C3").Value = 33.183996089999994
$ws.Range("D3").Value = 29.790700447499944
$ws.Range("E3").Value = 27.370778812500021

# Selection now only spans the edited columns
$ws.Range("B1:E3").Select()
